# Add KV Guardrails Assignment to Platform Management group (#1383)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rename the worksheet (propagates to the _FilterDatabase defined name) ---
$ws.Name = "ALZ Policy Assignments"

# --- Populate row 15 with the new "Enforce recommended guardrails for Azure
#     Key Vault" assignment (duplicate of the row 35 initiative, new release date) ---
$ws.Range("B15").Value = "Enforce recommendded guardrails for Azure Key Vault"
$ws.Range("C15").Value = "Enforce recommendded guardrails for Azure Key Vault"
$ws.Range("D15").Value = "Initiative"
$ws.Range("E15").Value = "Custom"
$ws.Range("F15").Value = "This initiative assignment enables recommended ALZ guardrails for Azure Key Vault."
$ws.Range("G15").Value = "Deny, Audit"
$ws.Range("H15").Value = "ENFORCE-GuardrailsKeyVaultPolicyAssignment.json"

# Add the AzAdvertizer hyperlink on I15 (address stored as the link's display
# text) and then set the cell's visible text independently, matching the
# existing pattern used by the other AzAdvertizer links in this sheet.
$null = $ws.Hyperlinks.Add($ws.Range("I15"), "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/Enforce-Guardrails-KeyVault.html", "", "", "https://www.azadvertizer.net/azpolicyinitiativesadvertizer/Enforce-Guardrails-KeyVault.html")
$ws.Range("I15").Value = "Enforce recommended guardrails for Azure Key Vault (azadvertizer.net)"

# Copy the formatting from the matching I35 cell (hyperlink style, no
# top-vertical alignment) onto I15 so it matches the other initiative rows.
$ws.Range("I35").Copy()
$ws.Range("I15").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("J15").Value = 45124

# Row 15 now wraps onto 3 lines like row 35 (same content) -> taller row.
$ws.Rows.Item(15).RowHeight = 43.2

# --- Update the sheet's saved view/selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("J15").Select()
